# Add a new "Docentes responsáveis" (responsible professor) entry.
# A new row is inserted above the existing "1033242 - Fábio Herbst
# Florenzano" row, pushing it (and everything below it) down by one row,
# and the new row is filled in with the additional professor's name in
# both column B and column C (mirroring the existing B/C duplication used
# throughout the sheet for the "current" / "modified (in red)" columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a whole new row at row 13 - this shifts the old row 13
# ("1033242 - Fábio Herbst Florenzano") down to row 14, and everything
# below it shifts down by one row as well (old row 24 -> new row 25).
$ws.Rows("13:13").Insert()

# The freshly inserted row pulls column A's bold-label formatting down
# into the new A13 even though that row has no label - clear it so no
# stray formatted-but-empty cell is left behind (matches rows like the
# existing 13/19/24 that only use columns B/C).
$ws.Range("A13").Clear()

# Fill in the new professor's name in both the "current" (B) and
# "modified" (C) columns, same as the existing name row below it.
$newTeacher = "5840963 - Daniela Camargo Vernilli"
$ws.Range("B13").Value = $newTeacher
$ws.Range("C13").Value = $newTeacher

# Copy the row-14 (the shifted-down original name row) formatting onto
# the new row 13 so B13/C13 pick up the same non-bold / wrap-text styles
# used for name cells, instead of the generic style Insert() applied.
$ws.Range("B14:C14").Copy()
$ws.Range("B13").PasteSpecial(-4122)
